# Insert a new row at position 241, shifting the existing rows 241-257 down
# to 242-258 (dimension grows from T257 to T258).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(241).Insert()

# Populate the new row 241 with the new pricing record.
$ws.Range("A241").Value = 5
$ws.Range("B241").Value = "Macroferia Regional de Talca"
$ws.Range("C241").Value = "Maule"
$ws.Range("D241").Value = 44714
$ws.Range("E241").Value = 7
$ws.Range("F241").Value = "Fruta"
$ws.Range("G241").Value = 100108
$ws.Range("H241").Value = "Tropicales y subtropicales"
$ws.Range("I241").Value = 100108005
$ws.Range("J241").Value = "Piña"
$ws.Range("K241").Value = "Caramelo"
$ws.Range("L241").Value = "Tercera"
$ws.Range("M241").Value = 150
$ws.Range("N241").Value = 18000
$ws.Range("O241").Value = 18000
$ws.Range("P241").Value = 18000
$ws.Range("Q241").Value = "`$/caja 16 unidades"
$ws.Range("R241").Value = "Ecuador"
$ws.Range("S241").Value = 1125
$ws.Range("T241").Value = 16
